$d = $word.ActiveDocument

$replacements = @(
    @{ old = "51×73="; new = "32×33=" },
    @{ old = "66×65="; new = "86×35=" },
    @{ old = "93×81="; new = "60×15=" },
    @{ old = "30×76="; new = "66×46=" },
    @{ old = "87×13="; new = "95×56=" },
    @{ old = "39×78="; new = "99×56=" },
    @{ old = "20×79="; new = "47×14=" },
    @{ old = "59×75="; new = "64×17=" },
    @{ old = "91×83="; new = "74×65=" },
    @{ old = "17×24="; new = "62×61=" },
    @{ old = "99×49="; new = "98×28=" },
    @{ old = "12×83="; new = "39×91=" },
    @{ old = "21×32="; new = "64×67=" },
    @{ old = "23×72="; new = "37×93=" },
    @{ old = "92×29="; new = "17×19=" },
    @{ old = "25×67="; new = "11×42=" },
    @{ old = "68×80="; new = "37×16=" },
    @{ old = "25×11="; new = "38×33=" },
    @{ old = "69×14="; new = "87×34=" },
    @{ old = "11×23="; new = "36×98=" },
    @{ old = "45×15="; new = "50×86=" },
    @{ old = "84×39="; new = "67×54=" },
    @{ old = "94×88="; new = "42×93=" },
    @{ old = "76×18="; new = "54×50=" },
    @{ old = "66×85="; new = "91×59=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
